# Updates the cryptos price/volume table with the latest scraped values.
# D column values that look numeric are apostrophe-prefixed so Excel
# retains them as literal text (matching the original formatting,
# e.g. preserving trailing zeros / multi-dot thousand separators).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'36.799.85"
$ws.Range("E2").Value = "  -1.26%  "
$ws.Range("D3").Value = "'1.991.64"
$ws.Range("E3").Value = "  -1.85%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "'255.99"
$ws.Range("E5").Value = "  +3.46%  "
$ws.Range("D6").Value = "'0.611"
$ws.Range("E6").Value = "  -2.07%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").Value = "'55.13"
$ws.Range("E8").Value = "  -7.91%  "
$ws.Range("E9").Value = "  -4.68%  "
$ws.Range("D10").Value = "'0.0763"
$ws.Range("E10").Value = "  -5.59%  "
$ws.Range("E11").Value = "  -3.05%  "
$ws.Range("B12").Value = "Chainlink"
$ws.Range("C12").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D12").Value = "'14.14"
$ws.Range("E12").Value = "  -7.10%  "
$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D13").Value = "'2.285.72"
$ws.Range("E13").Value = "  -1.62%  "
$ws.Range("D14").Value = "'21.28"
$ws.Range("E14").Value = "  -4.00%  "
$ws.Range("E15").Value = "  -7.12%  "
$ws.Range("E16").Value = "  -5.41%  "
$ws.Range("D17").Value = "'1.997.35"
$ws.Range("E17").Value = "  -1.54%  "
$ws.Range("D18").Value = "'36.712.75"
$ws.Range("E18").Value = "  -1.45%  "
$ws.Range("D19").Value = "'70.49"
$ws.Range("E19").Value = "  +0.26%  "
$ws.Range("D20").Value = "'0.0₃0821"
$ws.Range("E20").Value = "  -4.68%  "
$ws.Range("D21").Value = "'234.77"
$ws.Range("E21").Value = "  +1.93%  "
$ws.Range("E22").Value = "  -3.50%  "
$ws.Range("E23").Value = "  -0.14%  "
$ws.Range("E24").Value = "  -1.51%  "
$ws.Range("D25").Value = "'2.37"
$ws.Range("E25").Value = "  +0.57%  "
$ws.Range("D26").Value = "'163.80"
$ws.Range("E26").Value = "  -0.44%  "
$ws.Range("E27").Value = "  -5.50%  "
$ws.Range("D28").Value = "'19.35"
$ws.Range("E28").Value = "  -3.32%  "
$ws.Range("B29").Value = "ImmutableX"
$ws.Range("C29").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D29").Value = "'1.33"
$ws.Range("E29").Value = "  -3.92%  "
$ws.Range("B30").Value = "Kaspa"
$ws.Range("C30").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D30").Value = "'0.123"
$ws.Range("E30").Value = "  -9.53%  "
$ws.Range("E31").Value = "  -2.58%  "
$ws.Range("D32").Value = "'4.53"
$ws.Range("E32").Value = "  -4.93%  "
$ws.Range("D33").Value = "'0.0629"
$ws.Range("E33").Value = "  -6.31%  "
$ws.Range("D34").Value = "'4.34"
$ws.Range("E34").Value = "  -4.07%  "
$ws.Range("D35").Value = "'2.34"
$ws.Range("E35").Value = "  -9.26%  "
$ws.Range("D36").Value = "'3.48"
$ws.Range("E36").Value = "  -3.55%  "
$ws.Range("E37").Value = "  +0.65%  "
$ws.Range("E38").Value = "  +0.16%  "
$ws.Range("D39").Value = "'5.47"
$ws.Range("E39").Value = "  +1.36%  "
$ws.Range("D40").Value = "'3.00"
$ws.Range("E40").Value = "  +0.28%  "
$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").Value = "'1.17"
$ws.Range("E41").Value = "  -1.03%  "
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "'1.444.27"
$ws.Range("E42").Value = "  +4.82%  "
$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").Value = "'0.0208"
$ws.Range("E43").Value = "  -4.03%  "
$ws.Range("B44").Value = "Cronos"
$ws.Range("C44").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D44").Value = "'0.0910"
$ws.Range("E44").Value = "  -6.27%  "
$ws.Range("D45").Value = "'88.43"
$ws.Range("E45").Value = "  -3.65%  "
$ws.Range("D46").Value = "'15.48"
$ws.Range("E46").Value = "  -7.13%  "
$ws.Range("E47").Value = "  -4.02%  "
$ws.Range("E48").Value = "  -0.01%  "
$ws.Range("D49").Value = "'6.84"
$ws.Range("E49").Value = "  -9.58%  "
$ws.Range("D50").Value = "'2.178.15"
$ws.Range("E50").Value = "  -1.71%  "
$ws.Range("E51").Value = "  -9.44%  "